$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.139386177062988
$ws.Range("B1").Value = 2.116159439086914
$ws.Range("C1").Value = 10.13556289672852
$ws.Range("D1").Value = 2.524795532226562
$ws.Range("E1").Value = 1.286947965621948
